$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 233; this shifts all existing rows 233..280 down to 234..281
# and extends the sheet dimension to A1:R281, matching the target diff.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with the new weekly price-observation data.
$ws.Range("A233").Value = 7
$ws.Range("B233").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C233").Value = "Ñuble"
$ws.Range("D233").Value = 44511
$ws.Range("E233").Value = 16
$ws.Range("F233").Value = 100114014
$ws.Range("G233").Value = "Betarraga"
$ws.Range("H233").Value = "Sin especificar"
$ws.Range("I233").Value = "Primera"
$ws.Range("J233").Value = 300
$ws.Range("K233").Value = 700
$ws.Range("L233").Value = 800
$ws.Range("M233").Value = 750
$ws.Range("N233").Value = "$/paquete 5 unidades"
$ws.Range("O233").Value = "Región del Maule"
$ws.Range("P233").Value = 150
$ws.Range("Q233").Value = 5
$ws.Range("R233").Value = "Hortaliza"
